$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C20").Value = "'41"
$ws.Range("D20").Value = "'113500.00"
$ws.Range("C21").Value = "'117"
$ws.Range("D21").Value = "'342175.00"
$ws.Range("C23").Value = "'248"
$ws.Range("D23").Value = "'888553.74"
$ws.Range("C25").Value = "'7"
$ws.Range("D25").Value = "'20024.69"
$ws.Range("C27").Value = "'11"
$ws.Range("D27").Value = "'35500.00"
$ws.Range("C28").Value = "'12"
$ws.Range("D28").Value = "'33450.00"
$ws.Range("C30").Value = "'34"
$ws.Range("D30").Value = "'134557.00"
$ws.Range("C79").Value = "'377"
$ws.Range("D79").Value = "'1311082.04"
$ws.Range("C82").Value = "'13"
$ws.Range("D82").Value = "'38500.00"
$ws.Range("C83").Value = "'58"
$ws.Range("D83").Value = "'183200.00"
$ws.Range("C88").Value = "'90"
$ws.Range("D88").Value = "'200500.00"
$ws.Range("C119").Value = "'211"
$ws.Range("D119").Value = "'580500.00"
$ws.Range("C120").Value = "'48"
$ws.Range("D120").Value = "'136547.58"
$ws.Range("C121").Value = "'393"
$ws.Range("D121").Value = "'1551756.95"
$ws.Range("C125").Value = "'79"
$ws.Range("D125").Value = "'229243.68"
$ws.Range("C126").Value = "'35"
$ws.Range("D126").Value = "'142579.76"
$ws.Range("C129").Value = "'63"
$ws.Range("D129").Value = "'248773.75"
$ws.Range("C130").Value = "'105"
$ws.Range("D130").Value = "'265266.44"
$ws.Range("C137").Value = "'1513"
$ws.Range("D137").Value = "'3590218.81"
$ws.Range("C142").Value = "'812"
$ws.Range("D142").Value = "'2052940.00"
$ws.Range("C144").Value = "'285"
$ws.Range("D144").Value = "'691626.11"
$ws.Range("C148").Value = "'6"
$ws.Range("D148").Value = "'12000.00"
$ws.Range("C149").Value = "'47"
$ws.Range("D149").Value = "'124000.00"
$ws.Range("C150").Value = "'42"
$ws.Range("D150").Value = "'109500.00"
$ws.Range("C151").Value = "'116"
$ws.Range("D151").Value = "'301000.00"
$ws.Range("C152").Value = "'19"
$ws.Range("D152").Value = "'44000.00"
$ws.Range("C153").Value = "'74"
$ws.Range("D153").Value = "'199000.00"
$ws.Range("C156").Value = "'7"
$ws.Range("D156").Value = "'14000.00"
$ws.Range("C157").Value = "'46"
$ws.Range("D157").Value = "'113800.00"
$ws.Range("C158").Value = "'33"
$ws.Range("D158").Value = "'71984.00"
$ws.Range("C159").Value = "'10"
$ws.Range("D159").Value = "'20000.00"
$ws.Range("C160").Value = "'4"
$ws.Range("D160").Value = "'9500.00"
$ws.Range("C161").Value = "'13"
$ws.Range("D161").Value = "'38500.00"
$ws.Range("C162").Value = "'57"
$ws.Range("D162").Value = "'120000.00"
$ws.Range("C164").Value = "'19"
$ws.Range("D164").Value = "'63441.00"
$ws.Range("C165").Value = "'15"
$ws.Range("D165").Value = "'37500.00"
$ws.Range("C166").Value = "'46"
$ws.Range("D166").Value = "'124119.00"
$ws.Range("C167").Value = "'153"
$ws.Range("D167").Value = "'396000.00"
$ws.Range("C168").Value = "'11"
$ws.Range("D168").Value = "'22000.00"
$ws.Range("C169").Value = "'271"
$ws.Range("D169").Value = "'889608.27"
$ws.Range("C170").Value = "'14"
$ws.Range("D170").Value = "'50703.43"
$ws.Range("C172").Value = "'10"
$ws.Range("D172").Value = "'28000.00"
$ws.Range("C173").Value = "'47"
$ws.Range("D173").Value = "'143858.79"
$ws.Range("C174").Value = "'17"
$ws.Range("D174").Value = "'40000.00"
$ws.Range("C175").Value = "'26"
$ws.Range("D175").Value = "'65650.00"
$ws.Range("C176").Value = "'10"
$ws.Range("D176").Value = "'24500.00"
$ws.Range("C177").Value = "'46"
$ws.Range("D177").Value = "'168307.95"
$ws.Range("C178").Value = "'64"
$ws.Range("D178").Value = "'153000.00"
$ws.Range("C181").Value = "'107"
$ws.Range("D181").Value = "'281500.00"
$ws.Range("C182").Value = "'329"
$ws.Range("D182").Value = "'875788.00"
$ws.Range("C184").Value = "'574"
$ws.Range("D184").Value = "'2001174.27"
$ws.Range("C185").Value = "'21"
$ws.Range("D185").Value = "'75000.00"
$ws.Range("C189").Value = "'57"
$ws.Range("D189").Value = "'147926.00"
$ws.Range("C190").Value = "'67"
$ws.Range("D190").Value = "'158000.00"
$ws.Range("C192").Value = "'98"
$ws.Range("D192").Value = "'396004.50"
$ws.Range("C193").Value = "'125"
$ws.Range("D193").Value = "'269196.77"
